$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Step 1: duplicate the current (pre-update) rows 80-83 into new rows
# 84-87, preserving all existing values and the date number format.
# ------------------------------------------------------------------
$srcRows = 80, 81, 82, 83
$dstRows = 84, 85, 86, 87
$cols = "A", "B", "C", "D", "E", "F", "G", "H", "I", "J", "K", "L", "M", "N", "O", "P", "Q", "R", "S", "T"
$dateCols = "D"

for ($i = 0; $i -lt $srcRows.Length; $i++) {
    $src = $srcRows[$i]
    $dst = $dstRows[$i]
    foreach ($col in $cols) {
        $srcCell = $ws.Range("$col$src")
        $dstCell = $ws.Range("$col$dst")
        if ($dateCols -contains $col) {
            $dstCell.NumberFormat = $srcCell.NumberFormat()
        }
        $dstCell.Value = $srcCell.Value()
    }
}

# ------------------------------------------------------------------
# Step 2: update rows 80-83 with the new week's figures.
# ------------------------------------------------------------------

# Row 80 (Especial)
$ws.Range("D80").Value = 45142
$ws.Range("M80").Value = 170
$ws.Range("P80").Value = 7412
$ws.Range("S80").Value = 2471

# Row 81 (Primera)
$ws.Range("D81").Value = 45142
$ws.Range("M81").Value = 200

# Row 82 (Segunda)
$ws.Range("D82").Value = 45142
$ws.Range("M82").Value = 230
$ws.Range("N82").Value = 3500
$ws.Range("P82").Value = 3696
$ws.Range("S82").Value = 1232

# Row 83 (Tercera)
$ws.Range("D83").Value = 45142
$ws.Range("M83").Value = 230
$ws.Range("N83").Value = 2500
$ws.Range("P83").Value = 2826
$ws.Range("S83").Value = 942
